$wb = $excel.ActiveWorkbook

# Update both "展览" and "全部类型" sheets: F3 345 -> 347, F5 100 -> 105
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 347
    $ws.Range("F5").Value = 105
}
